$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "AddCustomer"
$ws.Range("A1").Value = "FirstName"
$ws.Range("B1").Value = "LastName"
$ws.Range("C1").Value = "PostCode"
$ws.Range("A2").Value = "Jack"
$ws.Range("B2").Value = "Daniels"
$ws.Range("C2").Value = "JD12345"
$ws.Range("A1:C2").EntireColumn.AutoFit()
$ws.Range("A2").Select()
